$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (pushes existing rows 2..29 down to 3..30)
$ws.Rows(2).Insert()

# Fill the new row 2 with the new award entry (Agency before Award title, to
# match the shared-string insertion order of the original edit)
$ws.Range("A2").Value2 = "Dr. Gunjan Mehta"
$ws.Range("C2").Value2 = "Microscopy Australia"
$ws.Range("B2").Value2 = "Travel Award from Microscopy Australia to attend the FoundingGIDE Imaging Data Ecosystem meeting in Brisbane, Australia, in October 2025,"

# D2 holds literal text ("2025 October") formatted like the other Month/Year
# text cells (e.g. D10/D24). Force text storage first so the date-like
# literal isn't auto-parsed into a date serial, then apply the real format.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "2025 October"
$ws.Range("D2").NumberFormat = $ws.Range("D10").NumberFormat

# Update the view: scroll back to show row 1 and select D3
$ws.Range("A1").Select()
$ws.Range("D3").Select()
